# Federated Learning presentation edit script
# - Renumber section titles on slides 9,10 (2 -> 3) and 11-14 (2. Quantum Machine
#   Learning -> 3. Quantum Federated Learning)
# - Slide 15: renumber/retitle (2. Quantum Machine Learning -> 4. Simulation) and
#   append a "Github" hyperlink mention after the existing subtitle text
# - Duplicate slide 15 to create a new slide 16 ("Related Works") with two GitHub
#   repository links, removing the picture shapes that came along with the dup.

$p = $ppt.ActivePresentation
$nbsp = [char]0x00A0

# --- Slide 9: "2. Quantum Federated Learning: " -> "3. Quantum Federated Learning: " (keep trailing NBSP)
$s9 = $p.Slides.Item(9)
$run = $s9.Shapes.Title.TextFrame.TextRange.Runs(1, 1)
$run.Text = "3. Quantum Federated Learning:" + $nbsp

# --- Slide 10: "2. Quantum Federated Learning: " -> "3. Quantum Federated Learning: " (keep trailing space)
$s10 = $p.Slides.Item(10)
$run = $s10.Shapes.Title.TextFrame.TextRange.Runs(1, 1)
$run.Text = "3. Quantum Federated Learning: "

# --- Slides 11-13: "2. Quantum Machine Learning: " -> "3. Quantum Federated Learning: " (keep trailing NBSP)
foreach ($idx in 11, 12, 13) {
    $s = $p.Slides.Item($idx)
    $run = $s.Shapes.Title.TextFrame.TextRange.Runs(1, 1)
    $run.Text = "3. Quantum Federated Learning:" + $nbsp
}

# --- Slide 14: same retitle
$s14 = $p.Slides.Item(14)
$run = $s14.Shapes.Title.TextFrame.TextRange.Runs(1, 1)
$run.Text = "3. Quantum Federated Learning:" + $nbsp

# --- Slide 15: "2. Quantum Machine Learning: " -> "4. Simulation: ", add trailing
# space to the "QFL-Deep Unfolding Network" run, and append a new "Github" run
# (hyperlinked-looking text run, matching the author's follow-up edit) after it.
$s15 = $p.Slides.Item(15)
$titleTr = $s15.Shapes.Title.TextFrame.TextRange
$run1 = $titleTr.Runs(1, 1)
$run1.Text = "4. Simulation:" + $nbsp
$run2 = $titleTr.Runs(2, 1)
$run2.Text = "QFL-Deep Unfolding Network "
$titleTr.InsertAfter("Github")

# --- Duplicate slide 15 -> becomes slide 16 at the end of the deck
$dup = $s15.Duplicate()
$s16 = $dup.Item(1)

# Remove the two picture shapes that were duplicated along with the slide
for ($i = $s16.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s16.Shapes.Item($i)
    if ($shp.Type -eq 13) {
        $shp.Delete()
    }
}

# Retitle slide 16
$s16Title = $s16.Shapes.Title
$s16TitleTr = $s16Title.TextFrame.TextRange
$s16TitleTr.Text = "Related Works"

# Replace the textbox body with the two GitHub project links
$tb = $s16.Shapes.Item("TextBox 8")
$tb.TextFrame.TextRange.Text = "https://github.com/WhiteByeBye/Post-Quantum-Secure-Blockchained-Federated-Learning.git`r" + "`r" + "`r" + "`rhttps://github.com/s222416822/PQC-QFL-Model.git`r"

$tbTr = $tb.TextFrame.TextRange
$link1 = $tbTr.Runs(1, 1)
$link1.ActionSettings(1).Hyperlink.Address = "https://github.com/WhiteByeBye/Post-Quantum-Secure-Blockchained-Federated-Learning.git"

$fullText = $tbTr.Text
$secondLinkText = "https://github.com/s222416822/PQC-QFL-Model.git"
$startPos = $fullText.IndexOf($secondLinkText) + 1
$link2 = $tbTr.Characters($startPos, $secondLinkText.Length)
$link2.ActionSettings(1).Hyperlink.Address = $secondLinkText
